# Update "想去人数" (want-to-go count) figures that changed between
# crawler runs (gh-pages data refresh, commit 456a3b4).
#
# Sheet "展览"   (展览信息 sheet, sheetId=1)
#   F2: 5571 -> 5585
#   F8: 368  -> 372
#
# Sheet "演出"   (演出信息 sheet, sheetId=2)
#   F2: 48 -> 49
#
# Sheet "全部类型" (aggregated sheet, sheetId=4)
#   F2: 5571 -> 5585
#   F8: 48   -> 49
#   F9: 368  -> 372

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5585
$wsExhibit.Range("F8").Value = 372

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 49

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5585
$wsAll.Range("F8").Value = 49
$wsAll.Range("F9").Value = 372
